$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename header cells to match new naming convention
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet right after "Monthly Trend", cloning its
# sheet-level properties (outline settings, page margins, formats) so the
# new sheet matches the workbook look-and-feel, then wipe its contents.
$wsMonthly.Copy([System.Type]::Missing, $wsMonthly)
$wsForecast = $wb.Worksheets.Item(3)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.Clear()

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds / PO_Forecast / yhat_lower / yhat_upper)
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 4
$wsForecast.Range("C2").Value = -71.52732936279892
$wsForecast.Range("D2").Value = 81.57901332119063
$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("B3").Value = 6
$wsForecast.Range("C3").Value = -69.65983246495878
$wsForecast.Range("D3").Value = 83.6661439097679
$wsForecast.Range("A4").Value = 45032.99999999999
$wsForecast.Range("B4").Value = 39
$wsForecast.Range("C4").Value = -39.27716789016991
$wsForecast.Range("D4").Value = 118.4146625590858
$wsForecast.Range("A5").Value = 45053.99999999999
$wsForecast.Range("B5").Value = 46
$wsForecast.Range("C5").Value = -33.68311015610893
$wsForecast.Range("D5").Value = 132.421190023514
$wsForecast.Range("A6").Value = 45060.99999999999
$wsForecast.Range("B6").Value = 49
$wsForecast.Range("C6").Value = -33.0437455797012
$wsForecast.Range("D6").Value = 130.4462576715427
$wsForecast.Range("A7").Value = 45067.99999999999
$wsForecast.Range("B7").Value = 52
$wsForecast.Range("C7").Value = -24.48763344372699
$wsForecast.Range("D7").Value = 123.6337591207963
$wsForecast.Range("A8").Value = 45074.99999999999
$wsForecast.Range("B8").Value = 54
$wsForecast.Range("C8").Value = -23.21352184387557
$wsForecast.Range("D8").Value = 131.4593001899384
$wsForecast.Range("A9").Value = 45081.99999999999
$wsForecast.Range("B9").Value = 57
$wsForecast.Range("C9").Value = -17.9932700435963
$wsForecast.Range("D9").Value = 130.4663422858428
$wsForecast.Range("A10").Value = 45088.99999999999
$wsForecast.Range("B10").Value = 59
$wsForecast.Range("C10").Value = -18.14776672689942
$wsForecast.Range("D10").Value = 136.6141204535682
$wsForecast.Range("A11").Value = 45095.99999999999
$wsForecast.Range("B11").Value = 62
$wsForecast.Range("C11").Value = -11.68031317542973
$wsForecast.Range("D11").Value = 143.5526410627554
$wsForecast.Range("A12").Value = 45102.99999999999
$wsForecast.Range("B12").Value = 64
$wsForecast.Range("C12").Value = -21.66772883788403
$wsForecast.Range("D12").Value = 135.6164550517646
$wsForecast.Range("A13").Value = 45109.99999999999
$wsForecast.Range("B13").Value = 67
$wsForecast.Range("C13").Value = -14.54379261267863
$wsForecast.Range("D13").Value = 150.2134779727839
$wsForecast.Range("A14").Value = 45116.99999999999
$wsForecast.Range("B14").Value = 69
$wsForecast.Range("C14").Value = -12.76100056703338
$wsForecast.Range("D14").Value = 142.936853169842
$wsForecast.Range("A15").Value = 45123.99999999999
$wsForecast.Range("B15").Value = 72
$wsForecast.Range("C15").Value = -7.148753782293134
$wsForecast.Range("D15").Value = 148.7693921095256
$wsForecast.Range("A16").Value = 45130.99999999999
$wsForecast.Range("B16").Value = 74
$wsForecast.Range("C16").Value = -7.481391958858327
$wsForecast.Range("D16").Value = 148.4988804990359
$wsForecast.Range("A17").Value = 45137.99999999999
$wsForecast.Range("B17").Value = 77
$wsForecast.Range("C17").Value = -2.118576653764385
$wsForecast.Range("D17").Value = 152.4322436699244
$wsForecast.Range("A18").Value = 45144.99999999999
$wsForecast.Range("B18").Value = 79
$wsForecast.Range("C18").Value = 1.76668045760222
$wsForecast.Range("D18").Value = 159.0997649799733
$wsForecast.Range("A19").Value = 45151.99999999999
$wsForecast.Range("B19").Value = 82
$wsForecast.Range("C19").Value = 0.2192060650645955
$wsForecast.Range("D19").Value = 159.7801900108355
$wsForecast.Range("A20").Value = 45158.99999999999
$wsForecast.Range("B20").Value = 84
$wsForecast.Range("C20").Value = 2.224567201983735
$wsForecast.Range("D20").Value = 159.8146617598035
$wsForecast.Range("A21").Value = 45165.99999999999
$wsForecast.Range("B21").Value = 87
$wsForecast.Range("C21").Value = 10.02860370561712
$wsForecast.Range("D21").Value = 171.2715769230483

# Apply the date style (copied from an existing date cell) to column A data rows
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A21").PasteSpecial(-4122)

# Restore the originally active sheet/selection (adding/copying a sheet shifts
# the active tab to the new sheet as a side effect)
$wsForecast.Range("A1").Select()
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()
